# "Started remote install story"
#
# The "Remotely run a script" story is generalised from EC2-specific
# wording to generic VM wording, the old "Specify AMI and size for EC2"
# story row is dropped (its row is taken over by the BDD-Framework rows
# that shift up), and a couple of now-unused trailing blank rows are
# trimmed from the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Remotely run a script" (theme "Remote") gets new,
#     VM-flavoured description/acceptance text (was EC2-specific). ---
$ws.Range("F8").Value = "Given An existing VM`nWhen I remotely run a script`nThen The script should be executed on the VM"
$ws.Range("E8").Value = "As a SysAdmin`nI want to remotely run a script on an existing VM`nSo that I can install packages on it "

# --- Row 11: used to hold the "Specify AMI and size for EC2" story;
#     that story is removed and the row now carries what used to be the
#     first BDD-Framework row's text. ---
$ws.Range("C11").Value = "BDD Framework"
$ws.Range("D11").Value = "Specify line number of feature file"
$ws.Range("E11").ClearContents()
$ws.Range("E11").WrapText = $false
$ws.Rows.Item(11).RowHeight = 41

# --- Row 12: now carries what used to be the second BDD-Framework row
#     (moved up), and picks up the "Done" marker. ---
$ws.Range("B12").Value = "H"
$ws.Range("C12").Value = "BDD Framework"
$ws.Range("D12").Value = "Specify feature file"
$ws.Range("I12").Value = "P"
$ws.Rows.Item(12).RowHeight = 43

# --- Row 13: emptied out now that its former content moved up to row 12. ---
$ws.Range("B13:D13").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Rows.Item(13).AutoFit()

# --- Row 40: drop the stray formatted-but-empty cells outside the
#     Backlog/InProgress/Done columns, matching the other blank rows below. ---
$ws.Range("B40:F40").Clear()
$ws.Range("J40:L40").Clear()

# --- Row 46 is no longer needed. ---
$ws.Rows.Item(46).Delete()

# --- Update the view: scroll/selection moved to E6. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("E6").Select()
